$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 989.4
$ws.Range("J17").Value = 905.8570999999999
$ws.Range("L17").Value = 2717.5713
$ws.Range("N17").Value = -3053.5713
$ws.Range("H33").Value = 254.18182
$ws.Range("I33").Value = 219.6
$ws.Range("J33").Value = 600
$ws.Range("K33").Value = 219.6
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = 9.400000000000006
$ws.Range("N33").Value = -1058
$ws.Range("H62").Value = 1079
$ws.Range("I62").Value = 998.75
$ws.Range("J62").Value = 1400
$ws.Range("K62").Value = 998.75
$ws.Range("L62").Value = 1400
$ws.Range("M62").Value = -374.75
$ws.Range("N62").Value = -2648
$ws.Range("H65").Value = 1079
$ws.Range("I65").Value = 998.75
$ws.Range("J65").Value = 1400
$ws.Range("K65").Value = 4993.75
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = -1873.75
$ws.Range("N65").Value = -13240
$ws.Range("H86").Value = 250634
$ws.Range("I86").Value = 933.3333
$ws.Range("J86").Value = 500334.66
$ws.Range("K86").Value = 933.3333
$ws.Range("L86").Value = 500334.66
$ws.Range("M86").Value = 189.6667
$ws.Range("N86").Value = -502580.66
$ws.Range("H89").Value = 250634
$ws.Range("I89").Value = 933.3333
$ws.Range("J89").Value = 500334.66
$ws.Range("K89").Value = 4666.6665
$ws.Range("L89").Value = 2501673.3
$ws.Range("M89").Value = 949.3334999999997
$ws.Range("N89").Value = -2512905.3
$ws.Range("H95").Value = 13532.333
$ws.Range("J95").Value = 13532.333
$ws.Range("L95").Value = 13532.333
$ws.Range("N95").Value = -19024.333
$ws.Range("H106").Value = 45477216
$ws.Range("I106").Value = 50020890
$ws.Range("K106").Value = 50020890
$ws.Range("M106").Value = -50020259
$ws.Range("H125").Value = 3149.75
$ws.Range("I125").Value = 3149.75
$ws.Range("K125").Value = 28347.75
$ws.Range("M125").Value = -25887.75
$ws.Range("H132").Value = 2941.24
$ws.Range("I132").Value = 1342.8636
$ws.Range("K132").Value = 4028.5908
$ws.Range("M132").Value = -1498.5908
$ws.Range("H135").Value = 1558.0667
$ws.Range("I135").Value = 1246.1
$ws.Range("J135").Value = 2182
$ws.Range("K135").Value = 11214.9
$ws.Range("L135").Value = 19638
$ws.Range("M135").Value = -8679.9
$ws.Range("N135").Value = -24708
$ws.Range("H137").Value = 1411
$ws.Range("I137").Value = 898.2857
$ws.Range("K137").Value = 2694.8571
$ws.Range("M137").Value = -144.8571000000002
$ws.Range("H141").Value = 5166.1665
$ws.Range("I141").Value = 3199.4
$ws.Range("K141").Value = 9598.200000000001
$ws.Range("M141").Value = -4418.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2803.6667
$ws.Range("I45").Value = 2155.5
$ws.Range("K45").Value = 2155.5
$ws.Range("M45").Value = -1778.5
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1314
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -6568
$ws.Range("N66").Value = -16864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3089823.2
$ws.Range("I105").Value = 4388769
$ws.Range("J105").Value = 4827.125
$ws.Range("K105").Value = 4388769
$ws.Range("L105").Value = 4827.125
$ws.Range("M105").Value = -4387022
$ws.Range("N105").Value = -8321.125
$ws.Range("H107").Value = 3249.5
$ws.Range("I107").Value = 3499
$ws.Range("K107").Value = 3499
$ws.Range("M107").Value = -1579

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2759.3333
$ws.Range("I31").Value = 1687.76
$ws.Range("K31").Value = 1687.76
$ws.Range("M31").Value = -1392.76
$ws.Range("H34").Value = 2759.3333
$ws.Range("I34").Value = 1687.76
$ws.Range("K34").Value = 1687.76
$ws.Range("M34").Value = -1485.76
$ws.Range("H62").Value = 103610.75
$ws.Range("I62").Value = 4814.6665
$ws.Range("J62").Value = 399999
$ws.Range("K62").Value = 4814.6665
$ws.Range("L62").Value = 399999
$ws.Range("M62").Value = -4190.6665
$ws.Range("N62").Value = -401247
$ws.Range("H65").Value = 103610.75
$ws.Range("I65").Value = 4814.6665
$ws.Range("J65").Value = 399999
$ws.Range("K65").Value = 24073.3325
$ws.Range("L65").Value = 1999995
$ws.Range("M65").Value = -20953.3325
$ws.Range("N65").Value = -2006235
$ws.Range("H96").Value = 2455.375
$ws.Range("J96").Value = 2455.375
$ws.Range("L96").Value = 2455.375
$ws.Range("N96").Value = -7947.375
$ws.Range("H99").Value = 10636.774
$ws.Range("I99").Value = 6421.4116
$ws.Range("J99").Value = 15755.429
$ws.Range("K99").Value = 6421.4116
$ws.Range("L99").Value = 15755.429
$ws.Range("M99").Value = -4923.4116
$ws.Range("N99").Value = -18751.429
$ws.Range("H109").Value = 64285.715
$ws.Range("J109").Value = 64285.715
$ws.Range("L109").Value = 64285.715
$ws.Range("N109").Value = -66365.715
$ws.Range("H126").Value = 10636.774
$ws.Range("I126").Value = 6421.4116
$ws.Range("J126").Value = 15755.429
$ws.Range("K126").Value = 19264.2348
$ws.Range("L126").Value = 47266.287
$ws.Range("M126").Value = -16794.2348
$ws.Range("N126").Value = -52206.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 345.63635
$ws.Range("J2").Value = 452.125
$ws.Range("L2").Value = 2712.75
$ws.Range("N2").Value = -2938.75
$ws.Range("H5").Value = 462.22223
$ws.Range("I5").Value = 462.22223
$ws.Range("K5").Value = 1386.66669
$ws.Range("M5").Value = -1274.66669
$ws.Range("H11").Value = 1015.8
$ws.Range("I11").Value = 726.6667
$ws.Range("K11").Value = 2180.0001
$ws.Range("M11").Value = -2040.0001
$ws.Range("H26").Value = 312.83334
$ws.Range("I26").Value = 144.5
$ws.Range("J26").Value = 397
$ws.Range("K26").Value = 433.5
$ws.Range("L26").Value = 1191
$ws.Range("M26").Value = -145.5
$ws.Range("N26").Value = -1767
$ws.Range("H40").Value = 47.153847
$ws.Range("I40").Value = 31.6
$ws.Range("K40").Value = 126.4
$ws.Range("M40").Value = -57.40000000000001
$ws.Range("H129").Value = 1671.4286
$ws.Range("I129").Value = 925
$ws.Range("J129").Value = 2666.6667
$ws.Range("K129").Value = 2775
$ws.Range("L129").Value = 8000.000100000001
$ws.Range("M129").Value = 2225
$ws.Range("N129").Value = -18000.0001
$ws.Range("H132").Value = 3433.9092
$ws.Range("I132").Value = 3109.75
$ws.Range("K132").Value = 27987.75
$ws.Range("M132").Value = -25457.75
$ws.Range("H135").Value = 462.22223
$ws.Range("I135").Value = 462.22223
$ws.Range("K135").Value = 4160.00007
$ws.Range("M135").Value = -1625.00007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 1700
$ws.Range("I33").Value = 1700
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1700
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1448
$ws.Range("N33").ClearContents()
$ws.Range("H102").Value = 1435.25
$ws.Range("I102").Value = 1514.909
$ws.Range("K102").Value = 1514.909
$ws.Range("M102").Value = 107.0909999999999
$ws.Range("H122").Value = 50696.57
$ws.Range("I122").Value = 2570.2144
$ws.Range("J122").Value = 146949.28
$ws.Range("K122").Value = 7710.6432
$ws.Range("L122").Value = 440847.84
$ws.Range("M122").Value = -5260.6432
$ws.Range("N122").Value = -445747.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1218.2
$ws.Range("I16").Value = 1197.75
$ws.Range("K16").Value = 1197.75
$ws.Range("M16").Value = -1027.75
$ws.Range("H46").Value = 3206.5833
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 3354.1428
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 3354.1428
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -3730.1428
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 7666.3335
$ws.Range("I100").Value = 7000
$ws.Range("K100").Value = 7000
$ws.Range("M100").Value = -6459
$ws.Range("H132").Value = 169033.83
$ws.Range("I132").Value = 202400.8
$ws.Range("J132").Value = 2199
$ws.Range("K132").Value = 607202.3999999999
$ws.Range("L132").Value = 6597
$ws.Range("M132").Value = -604672.3999999999
$ws.Range("N132").Value = -11657

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 24974.5
$ws.Range("J69").Value = 24974.5
$ws.Range("L69").Value = 24974.5
$ws.Range("N69").Value = -26472.5
$ws.Range("H72").Value = 24974.5
$ws.Range("J72").Value = 24974.5
$ws.Range("L72").Value = 74923.5
$ws.Range("N72").Value = -82411.5
$ws.Range("H124").Value = 69999.5
$ws.Range("J124").Value = 69999.5
$ws.Range("L124").Value = 69999.5
$ws.Range("N124").Value = -79819.5
$ws.Range("H125").Value = 96297
$ws.Range("J125").Value = 96297
$ws.Range("L125").Value = 96297
$ws.Range("N125").Value = -106137
